# Auto-generated edit script applying value changes per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 8800.799999999999
$ws.Range("I106").Value = 8499.5
$ws.Range("K106").Value = 8499.5
$ws.Range("M106").Value = -7868.5

$ws.Range("H107").Value = 814.0769
$ws.Range("I107").Value = 886.8889
$ws.Range("K107").Value = 886.8889
$ws.Range("M107").Value = 1033.1111

$ws.Range("H132").Value = 9668.559999999999
$ws.Range("I132").Value = 9726.695
$ws.Range("K132").Value = 29180.085
$ws.Range("M132").Value = -26650.085

$ws.Range("H138").Value = 1397
$ws.Range("I138").Value = 1397
$ws.Range("K138").Value = 4191
$ws.Range("M138").Value = 949

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3928.7646
$ws.Range("I132").Value = 4202.4287
$ws.Range("K132").Value = 12607.2861
$ws.Range("M132").Value = -10077.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 223.72728
$ws.Range("I22").Value = 237.2
$ws.Range("K22").Value = 237.2
$ws.Range("M22").Value = -64.19999999999999

$ws.Range("H107").Value = 9366.916999999999
$ws.Range("I107").Value = 9058.143
$ws.Range("K107").Value = 9058.143
$ws.Range("M107").Value = -7138.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 550380.5
$ws.Range("I16").Value = 1100011
$ws.Range("K16").Value = 1100011
$ws.Range("M16").Value = -1099724

$ws.Range("H22").Value = 1111.8572
$ws.Range("J22").Value = 1271.8334
$ws.Range("L22").Value = 1271.8334
$ws.Range("N22").Value = -1971.8334

$ws.Range("H58").Value = 3203.3
$ws.Range("I58").Value = 1755.375
$ws.Range("J58").Value = 8995
$ws.Range("K58").Value = 1755.375
$ws.Range("L58").Value = 8995
$ws.Range("M58").Value = -1552.375
$ws.Range("N58").Value = -9401

$ws.Range("H113").Value = 550380.5
$ws.Range("I113").Value = 1100011
$ws.Range("K113").Value = 1100011
$ws.Range("M113").Value = -1097841

$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws.Range("H134").Value = 2648.75
$ws.Range("I134").Value = 2648.75
$ws.Range("K134").Value = 7946.25
$ws.Range("M134").Value = -5411.25

$ws.Range("H136").Value = 3203.3
$ws.Range("I136").Value = 1755.375
$ws.Range("J136").Value = 8995
$ws.Range("K136").Value = 5266.125
$ws.Range("L136").Value = 26985
$ws.Range("M136").Value = -2716.125
$ws.Range("N136").Value = -32085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 117263.89
$ws.Range("I4").Value = 117691.18
$ws.Range("K4").Value = 353073.54
$ws.Range("M4").Value = -352961.54

$ws.Range("H23").Value = 144.83333
$ws.Range("I23").Value = 55
$ws.Range("K23").Value = 165
$ws.Range("M23").Value = 70

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H103").Value = 1484.3334
$ws.Range("J103").Value = 1484.3334
$ws.Range("L103").Value = 4453.0002
$ws.Range("N103").Value = -6211.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 1000
$ws.Range("M70").Value = -730

$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 1000
$ws.Range("M73").Value = -64

$ws.Range("H80").Value = 4498.8335
$ws.Range("I80").Value = 4398.8
$ws.Range("J80").Value = 4999
$ws.Range("K80").Value = 4398.8
$ws.Range("L80").Value = 4999
$ws.Range("M80").Value = -3400.8
$ws.Range("N80").Value = -6995

$ws.Range("H83").Value = 4498.8335
$ws.Range("I83").Value = 4398.8
$ws.Range("J83").Value = 4999
$ws.Range("K83").Value = 21994
$ws.Range("L83").Value = 24995
$ws.Range("M83").Value = -17002
$ws.Range("N83").Value = -34979

$ws.Range("H102").Value = 2015.9375
$ws.Range("I102").Value = 1805
$ws.Range("K102").Value = 1805
$ws.Range("M102").Value = -183

$ws.Range("H113").Value = 2568.5715
$ws.Range("I113").Value = 1287.2727
$ws.Range("K113").Value = 1287.2727
$ws.Range("M113").Value = 882.7273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7410.75
$ws.Range("I7").Value = 7410.75
$ws.Range("K7").Value = 7410.75
$ws.Range("M7").Value = -7298.75

$ws.Range("H16").Value = 475
$ws.Range("I16").Value = 475
$ws.Range("K16").Value = 475
$ws.Range("M16").Value = -305

$ws.Range("H22").Value = 1571.6666
$ws.Range("J22").Value = 816.3333
$ws.Range("L22").Value = 816.3333
$ws.Range("N22").Value = -1406.3333

$ws.Range("H27").Value = 1571.6666
$ws.Range("J27").Value = 816.3333
$ws.Range("L27").Value = 816.3333
$ws.Range("N27").Value = -1030.3333

$ws.Range("H55").Value = 1351.8334
$ws.Range("I55").Value = 1696.4
$ws.Range("J55").Value = 1105.7142
$ws.Range("K55").Value = 1696.4
$ws.Range("L55").Value = 1105.7142
$ws.Range("M55").Value = -1523.4
$ws.Range("N55").Value = -1451.7142

$ws.Range("H126").Value = 7410.75
$ws.Range("I126").Value = 7410.75
$ws.Range("K126").Value = 22232.25
$ws.Range("M126").Value = -19762.25

$ws.Range("H136").Value = 1666.6666
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 40066.668
$ws.Range("I2").Value = 40066.668
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 40066.668
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -39954.668
$ws.Range("N2").ClearContents()

$ws.Range("H81").Value = 981.6667
$ws.Range("I81").Value = 900
$ws.Range("J81").Value = 998
$ws.Range("K81").Value = 1800
$ws.Range("L81").Value = 1996
$ws.Range("M81").Value = -739
$ws.Range("N81").Value = -4118

$ws.Range("H84").Value = 981.6667
$ws.Range("I84").Value = 900
$ws.Range("J84").Value = 998
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 9980
$ws.Range("M84").Value = -3696
$ws.Range("N84").Value = -20588

$ws.Range("H96").Value = 1350
$ws.Range("I96").Value = 1350
$ws.Range("K96").Value = 1350
$ws.Range("M96").Value = 23

$ws.Range("H107").Value = 459.5
$ws.Range("I107").Value = 451.4
$ws.Range("K107").Value = 1354.2
$ws.Range("M107").Value = 565.8000000000002

$ws.Range("H132").Value = 1975.2222
$ws.Range("I132").Value = 1975.2222
$ws.Range("K132").Value = 5925.6666
$ws.Range("M132").Value = -3395.6666
